$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '321.29'
Set-TextValue 'E2' '6.21%'
Set-TextValue 'D3' '49.09'
Set-TextValue 'E3' '11.19%'
Set-TextValue 'D4' '5.321'
Set-TextValue 'E4' '4.55%'
Set-TextValue 'D5' '0.08061'
Set-TextValue 'E5' '4.60%'
Set-TextValue 'D6' '4.601'
Set-TextValue 'E6' '4.15%'
Set-TextValue 'D7' '1.345'
Set-TextValue 'E7' '28.47%'
Set-TextValue 'D8' '1.643'
Set-TextValue 'E8' '1.45%'
Set-TextValue 'D9' '0.1282'
Set-TextValue 'E9' '0.80%'
Set-TextValue 'E10' '5.59%'
Set-TextValue 'D11' '0.09631'
Set-TextValue 'E11' '5.73%'
Set-TextValue 'D12' '0.04719'
Set-TextValue 'E12' '13.47%'
Set-TextValue 'D13' '0.1048'
Set-TextValue 'E13' '0.13%'
Set-TextValue 'D14' '0.001324'
Set-TextValue 'E14' '3.78%'
Set-TextValue 'D15' '0.04193'
Set-TextValue 'E15' '0.12%'
Set-TextValue 'D16' '0.005895'
Set-TextValue 'E16' '2.21%'
Set-TextValue 'D17' '3.344'
Set-TextValue 'E17' '-0.03%'
Set-TextValue 'D18' '2.446'
Set-TextValue 'E18' '4.95%'
Set-TextValue 'E19' '4.93%'
Set-TextValue 'D20' '8.031'
Set-TextValue 'E20' '-0.86%'
Set-TextValue 'D21' '0.1365'
Set-TextValue 'E21' '-1.62%'
Set-TextValue 'D22' '0.3091'
Set-TextValue 'E22' '-2.72%'
Set-TextValue 'D23' '0.001314'
Set-TextValue 'E23' '2.43%'
Set-TextValue 'D24' '0.004336'
Set-TextValue 'E24' '-1.90%'
Set-TextValue 'D25' '0.0001347'
Set-TextValue 'E25' '-0.09%'
Set-TextValue 'D26' '0.0003537'
Set-TextValue 'D38' '0.02721'
Set-TextValue 'E38' '8.81%'
Set-TextValue 'D39' '0.06000'
Set-TextValue 'E39' '13.18%'
Set-TextValue 'E40' '82.98%'
Set-TextValue 'D41' '0.008022'
Set-TextValue 'E41' '4.10%'
Set-TextValue 'D42' '0.1467'
Set-TextValue 'E42' '8.56%'
Set-TextValue 'D43' '0.007899'
Set-TextValue 'E43' '7.53%'
Set-TextValue 'D44' '0.008643'
Set-TextValue 'E44' '14.69%'
Set-TextValue 'D45' '0.3503'
Set-TextValue 'E45' '16.28%'
Set-TextValue 'D46' '0.00006970'
Set-TextValue 'E46' '4.45%'
Set-TextValue 'E47' '-0.11%'
Set-TextValue 'D48' '0.05962'
Set-TextValue 'E48' '38.38%'
Set-TextValue 'D49' '0.003997'
Set-TextValue 'E49' '-4.85%'
Set-TextValue 'E50' '-0.11%'
Set-TextValue 'D51' '0.0001997'
Set-TextValue 'E51' '-0.11%'
